$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.251.97"
$ws.Range("E2").Value = "  -1.42%  "

# Row 3
$ws.Range("D3").Value = "2.929.75"
$ws.Range("E3").Value = "  -2.94%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.32"
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.00%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("E8").Value = "  -0.71%  "

# Row 9
$ws.Range("D9").Value = "2.924.12"
$ws.Range("E9").Value = "  -3.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.81%  "

# Row 11
$ws.Range("E11").Value = "  -4.14%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.35%  "

# Row 13
$ws.Range("E13").Value = "  -3.77%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.25%  "

# Row 15
$ws.Range("E15").Value = "  -0.82%  "

# Row 16
$ws.Range("D16").Value = "65.240.06"
$ws.Range("E16").Value = "  -1.40%  "

# Row 17
$ws.Range("D17").Value = "3.417.89"
$ws.Range("E17").Value = "  -2.81%  "

# Row 18
$ws.Range("E18").Value = "  +0.20%  "

# Row 19
$ws.Range("D19").Value = "2.928.83"
$ws.Range("E19").Value = "  -2.32%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +12.23%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "442.85"
$ws.Range("D21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.690"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.67%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.86%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.92%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.60%  "

# Row 27
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.38%  "

# Row 28
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "

# Row 30
$ws.Range("E30").Value = "  -0.45%  "

# Row 31
$ws.Range("E31").Value = "  -1.66%  "

# Row 32
$ws.Range("E32").Value = "  -4.74%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.42%  "

# Row 35
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.971"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.73%  "

# Row 37
$ws.Range("E37").Value = "  -1.69%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.24%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "44.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.68%  "

# Row 40
$ws.Range("E40").Value = "  -9.32%  "

# Row 41
$ws.Range("E41").Value = "  -1.46%  "

# Row 42
$ws.Range("E42").Value = "  -2.57%  "

# Row 43
$ws.Range("E43").Value = "  -7.86%  "

# Row 44
$ws.Range("E44").Value = "  +0.25%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "381.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.22%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0350"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.32%  "

# Row 47
$ws.Range("D47").Value = "2.695.85"
$ws.Range("E47").Value = "  -3.59%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.20%  "

# Row 50
$ws.Range("E50").Value = "  +4.43%  "

# Row 51
$ws.Range("E51").Value = "  +0.01%  "
